$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 1; 3 = 2; 4 = 1; 5 = 2; 6 = 0; 7 = 1; 8 = 2; 9 = 2; 10 = 2;
    11 = 1; 12 = 0; 13 = 2; 14 = 1; 15 = 2; 16 = 3; 17 = 3; 18 = 0; 19 = 1; 20 = 1;
    21 = 1; 22 = 3; 23 = 3; 24 = 2; 25 = 2; 26 = 0; 27 = 0; 28 = 1; 29 = 3; 30 = 1;
    32 = 2; 33 = 3; 34 = 0; 35 = 2; 36 = 2; 37 = 1; 38 = 1; 39 = 3;
    40 = 0; 41 = 0; 42 = 2; 43 = 2; 44 = 2; 45 = 1; 46 = 1; 47 = 2; 48 = 2; 49 = 0;
    50 = 2; 51 = 2; 52 = 0; 53 = 2; 54 = 2; 55 = 2; 56 = 2; 57 = 0; 58 = 1; 59 = 1;
    60 = 3; 61 = 2; 62 = 0; 63 = 1; 64 = 2; 65 = 2;
    67 = 0; 68 = 2; 69 = 1; 70 = 0; 71 = 0; 72 = 2; 73 = 1; 74 = 1; 75 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
